# Apply edit: insert a new row for "climate_change_factor_gnrl_hydropower_availability"
# above the "elasticity_gnrl_rate_occupancy_to_gdppc" row (row 4) on sheet "strategy_id-0",
# shifting that row and all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Insert a new row at row 4; existing rows 4-11 shift down to 5-12.
$ws.Rows.Item(4).Insert()

# Columns H..AS as they appear in row 1 (35 data columns total: H=max_35, I=min_35, J..AS = years 0..35)
$cols = @("H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
          "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS")

# Populate the new row 4 (mirrors the metadata columns of the surrounding rows;
# C4:G4 are left blank, same as the empty placeholder cells in the neighboring rows)
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"

$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

foreach ($col in $cols[2..($cols.Length-1)]) {
    $ws.Range("$col`4").Value = 1
}

$wb.Save()
